$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: fill in the missing end time (D19) and duration (E19) ---
$ws.Range("C19").Copy()
$ws.Range("D19").PasteSpecial(-4122)          # xlPasteFormats: reuse the time style
$ws.Range("D19").Value = 0.72916666666666663  # 17:30

$ws.Range("C19").Copy()
$ws.Range("E19").PasteSpecial(-4122)          # xlPasteFormats: reuse the time style
$ws.Range("E19").Formula = "=D19-C19"

# --- Row 20: new "Demo" time-logging entry ---
$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)          # xlPasteFormats: reuse the date style
$ws.Range("B20").Value = 45999                # 2025-12-08

$ws.Range("C19").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 0.72916666666666663  # 17:30

$ws.Range("C19").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 0.79166666666666663  # 19:00

$ws.Range("C19").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Formula = "=D20-C20"

$ws.Range("F20").Value = "Demo"

$excel.CutCopyMode = 0

# Matches the cursor resting on the newly entered end-time cell
$ws.Range("D20").Select()
